{"js": "// Replace the merge-field placeholder text \"+++=prixintra_calc+++\" with\n// \"+++=prixintra+++\" (dropping the \"_calc\" suffix while keeping the\n// leading \"+++=prixintra\" and trailing \"+++\" markers, now as a single run).\nconst body = context.document.body;\n\nconst results = body.search(\"+++=prixintra_calc+++\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text '+++=prixintra_calc+++' not found in document body.\");\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"+++=prixintra+++\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Collapse the merge-field placeholder \"+++=prixintra_calc+++\" down to\n# \"+++=prixintra+++\" by deleting the \"_calc\" substring (it was split across\n# its own run in the original markup; removing it lets Word merge the\n# surrounding text back into a single run).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"_calc\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"\"\n$find.Forward = $true\n$find.Wrap = 1          # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
